$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 2025-09 row (row 22) stats
$ws.Range("B22").Value = 6303
$ws.Range("C22").Value = 997
$ws.Range("D22").Value = 5866735
$ws.Range("E22").Value = 930.7845470410915
$ws.Range("F22").Value = 8.50404544672061
$ws.Range("G22").Value = 4.288702928870292
$ws.Range("H22").Value = 27.58320357692743
